$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace placeholder staff names (column A, rows 4-33) with real employee names
$ws.Range("A4").Value = "川畑　拓也"
$ws.Range("A5").Value = "中田　拓実"
$ws.Range("A6").Value = "高橋　真央"
$ws.Range("A7").Value = "友安　美琴"
$ws.Range("A8").Value = "降矢　悠司"
$ws.Range("A9").Value = "横山　加奈"
$ws.Range("A10").Value = "渡邊　夏希"
$ws.Range("A11").Value = "髙田　実歩"
$ws.Range("A12").Value = "藤本　春佳"
$ws.Range("A13").Value = "永田　慎治"
$ws.Range("A14").Value = "花房　秀明"
$ws.Range("A15").Value = "山田　伸明"
$ws.Range("A16").Value = "山上　恵実"
$ws.Range("A17").Value = "徳重　駿介"
$ws.Range("A18").Value = "水越　智哉"
$ws.Range("A19").Value = "杉村　青空"
$ws.Range("A20").Value = "長坂　有里子"
$ws.Range("A21").Value = "藤本　望"
$ws.Range("A22").Value = "安達　由佳"
$ws.Range("A23").Value = "西尾　貴宏"
$ws.Range("A24").Value = "木口　瞳子"
$ws.Range("A25").Value = "上道　啓太"
$ws.Range("A26").Value = "松田　瑠美"
$ws.Range("A27").Value = "水野　真吾"
$ws.Range("A28").Value = "石嶺　利樹"
$ws.Range("A29").Value = "末堂　厚"
$ws.Range("A30").Value = "加藤　清澄"
$ws.Range("A31").Value = "鎬昂　昇"
$ws.Range("A32").Value = "本部　以蔵"
$ws.Range("A33").Value = "松本　梢江"

# Highlight the last-month attribute columns (C:G) for each staff row in red,
# matching the existing center/middle alignment already applied to that range.
$ws.Range("C4:G33").Font.ColorIndex = 3

